$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (STAFF_ID) becomes text-formatted "110".."113" instead of numeric values
$eRange = $ws.Range("E12:E16")
$eRange.NumberFormat = "@"
$ws.Range("E12").Value = "110"
$ws.Range("E13").Value = "111"
$ws.Range("E14").Value = "112"
$ws.Range("E15").Value = "113"
$ws.Range("E16").Value = "110"

# Rename subject/class id "SC11" -> "SCI11" wherever it appears (rows 14 and 16)
$ws.Range("C14").Value = "SCI11"
$ws.Range("C16").Value = "SCI11"

# Update the active selection to match the saved view state
$ws.Range("E17").Select()
